$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price and volume columns for rows with changed values
$ws.Range('D2').Value = '34.081.65'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '1.790.17'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '227.00'
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '32.31'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +4.00%  '
$ws.Range('D10').Value = '0.0688'
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = '2.047.03'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').Value = '11.38'
$ws.Range('D14').Value = '1.794.64'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').Value = '0.623'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').Value = '34.082.27'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').Value = '4.19'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = '68.07'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = '243.76'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '0.0₃0782'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').Value = '10.94'
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('D23').Value = '4.10'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('E24').Value = '  -2.71%  '
$ws.Range('D25').Value = '161.92'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').Value = '7.21'
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('D27').Value = '16.29'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').Value = '1.24'
$ws.Range('E30').Value = '  +2.37%  '
$ws.Range('D31').Value = '0.0522'
$ws.Range('E31').Value = '  +1.90%  '
$ws.Range('D32').Value = '3.66'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').Value = '3.61'
$ws.Range('E33').Value = '  +3.43%  '
$ws.Range('E34').Value = '  +1.50%  '
$ws.Range('D35').Value = '1.413.82'
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('E37').Value = '  +2.83%  '
$ws.Range('E38').Value = '  +7.88%  '
$ws.Range('E39').Value = '  -0.57%  '
$ws.Range('D40').Value = '80.76'
$ws.Range('E40').Value = '  +3.30%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').Value = '0.922'
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '13.40'
$ws.Range('E44').Value = '  +8.91%  '
$ws.Range('D45').Value = '0.0₆0138'
$ws.Range('E45').Value = '  -4.05%  '
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('D49').Value = '107.18'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').Value = '1.947.99'
$ws.Range('E50').Value = '  +0.43%  '

# Rows 46 and 47 swapped coin identity (Kaspa <-> FraxShare) with new values
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '6.05'
$ws.Range('E46').Value = '  +3.25%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').Value = '0.0507'
$ws.Range('E47').Value = '  +2.15%  '
